$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.491.11'
$ws.Cells.Item(2, 5).Value = '  -3.54%  '
$ws.Cells.Item(3, 4).Value = '1.993.00'
$ws.Cells.Item(3, 5).Value = '  -6.34%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.009'
$ws.Cells.Item(4, 5).Value = '  +0.26%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '329.12'
$ws.Cells.Item(5, 5).Value = '  -5.35%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.008'
$ws.Cells.Item(6, 5).Value = '  +0.26%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5001'
$ws.Cells.Item(7, 5).Value = '  -4.72%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.4212'
$ws.Cells.Item(8, 5).Value = '  -6.37%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '52.42'
$ws.Cells.Item(9, 5).Value = '  -3.18%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.08881'
$ws.Cells.Item(10, 5).Value = '  -5.35%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.119'
$ws.Cells.Item(11, 5).Value = '  -5.54%  '
$ws.Cells.Item(12, 2).Value = 'WrappedEther'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(12, 4).Value = '2.110.36'
$ws.Cells.Item(12, 5).Value = '  +2.00%  '
$ws.Cells.Item(13, 2).Value = 'Solana'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '23.32'
$ws.Cells.Item(13, 5).Value = '  -8.33%  '
$ws.Cells.Item(14, 2).Value = 'Chainlink'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '8.088'
$ws.Cells.Item(14, 5).Value = '  -7.53%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '6.503'
$ws.Cells.Item(15, 5).Value = '  -6.77%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '96.07'
$ws.Cells.Item(16, 5).Value = '  -7.06%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '1.010'
$ws.Cells.Item(17, 5).Value = '  +0.27%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.00001106'
$ws.Cells.Item(18, 5).Value = '  -5.63%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06625'
$ws.Cells.Item(19, 5).Value = '  -1.53%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '19.72'
$ws.Cells.Item(20, 5).Value = '  -8.58%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.008'
$ws.Cells.Item(21, 5).Value = '  +0.33%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.958'
$ws.Cells.Item(22, 5).Value = '  -5.93%  '
$ws.Cells.Item(23, 4).Value = '29.520.78'
$ws.Cells.Item(23, 5).Value = '  -3.37%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.86'
$ws.Cells.Item(24, 5).Value = '  -7.32%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.279'
$ws.Cells.Item(25, 5).Value = '  -2.43%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '158.09'
$ws.Cells.Item(26, 5).Value = '  -3.61%  '
$ws.Cells.Item(27, 2).Value = 'EthereumClassic'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '20.61'
$ws.Cells.Item(27, 5).Value = '  -7.61%  '
$ws.Cells.Item(28, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '6.559'
$ws.Cells.Item(28, 5).Value = '  -6.37%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.330'
$ws.Cells.Item(29, 5).Value = '  -8.50%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '127.67'
$ws.Cells.Item(30, 5).Value = '  -5.39%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.050'
$ws.Cells.Item(31, 5).Value = '  -9.88%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.09916'
$ws.Cells.Item(32, 5).Value = '  -6.64%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.555'
$ws.Cells.Item(33, 5).Value = '  -13.34%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.825'
$ws.Cells.Item(34, 5).Value = '  -7.79%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.792'
$ws.Cells.Item(35, 5).Value = '  -4.33%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '9.571'
$ws.Cells.Item(36, 5).Value = '  -10.52%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.02459'
$ws.Cells.Item(37, 5).Value = '  -7.82%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.06361'
$ws.Cells.Item(38, 5).Value = '  -7.65%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.287'
$ws.Cells.Item(39, 5).Value = '  -3.67%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.6513'
$ws.Cells.Item(40, 5).Value = '  -8.91%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '11.68'
$ws.Cells.Item(41, 5).Value = '  -8.45%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.2063'
$ws.Cells.Item(42, 5).Value = '  -8.77%  '
$ws.Cells.Item(43, 5).Value = '  +0.32%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.6333'
$ws.Cells.Item(44, 5).Value = '  -9.14%  '
$ws.Cells.Item(45, 2).Value = 'NEARProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.207'
$ws.Cells.Item(45, 5).Value = '  -7.96%  '
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '13.37'
$ws.Cells.Item(46, 5).Value = '  -9.27%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.271'
$ws.Cells.Item(47, 5).Value = '  -0.72%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '3.512'
$ws.Cells.Item(48, 5).Value = '  -3.41%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.00000000333'
$ws.Cells.Item(49, 5).Value = '  -4.40%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.07005'
$ws.Cells.Item(50, 5).Value = '  -3.21%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.140'
$ws.Cells.Item(51, 5).Value = '  -5.59%  '
